$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 2).Value = 0.9903733673316458
$ws.Cells.Item(1, 3).Value = -0.0002537796048658035
$ws.Cells.Item(1, 4).Value = -1.40006462893066
$ws.Cells.Item(1, 5).Value = 0.1709853632432805
$ws.Cells.Item(1, 6).Value = 1.570796390562869
$ws.Cells.Item(2, 2).Value = 0.9921109396706966
$ws.Cells.Item(2, 3).Value = -0.0002475824313470903
$ws.Cells.Item(2, 4).Value = -1.400018602685125
$ws.Cells.Item(2, 5).Value = 0.1710251922030726
$ws.Cells.Item(2, 6).Value = 1.570796390573616
$ws.Cells.Item(3, 2).Value = 1.002947235608045
$ws.Cells.Item(3, 3).Value = -0.000208934022463834
$ws.Cells.Item(3, 4).Value = -1.39973156196383
$ws.Cells.Item(3, 5).Value = 0.1712735838156022
$ws.Cells.Item(3, 6).Value = 1.570796390640639
$ws.Cells.Item(4, 2).Value = 1.02857456005494
$ws.Cells.Item(4, 3).Value = -0.0001175323738254093
$ws.Cells.Item(4, 4).Value = -1.399052724311498
$ws.Cells.Item(4, 5).Value = 0.1718610181641091
$ws.Cells.Item(4, 6).Value = 1.570796390799146
$ws.Cells.Item(5, 2).Value = 1.07145808442964
$ws.Cells.Item(5, 3).Value = 0.00003541472119404524
$ws.Cells.Item(5, 4).Value = -1.397916790239622
$ws.Cells.Item(5, 5).Value = 0.1728440023712562
$ws.Cells.Item(5, 6).Value = 1.570796391064382
$ws.Cells.Item(6, 2).Value = 1.131373702239576
$ws.Cells.Item(6, 3).Value = 0.0002491079710830877
$ws.Cells.Item(6, 4).Value = -1.396329696065339
$ws.Cells.Item(6, 5).Value = 0.1742173994258922
$ws.Cells.Item(6, 6).Value = 1.570796391434963
$ws.Cells.Item(7, 2).Value = 1.205945884663518
$ws.Cells.Item(7, 3).Value = 0.0005150748858204336
$ws.Cells.Item(7, 4).Value = -1.394354366750301
$ws.Cells.Item(7, 5).Value = 0.1759267570098149
$ws.Cells.Item(7, 6).Value = 1.570796391896196
$ws.Cells.Item(8, 2).Value = 1.291185536133739
$ws.Cells.Item(8, 3).Value = 0.0008190880765024157
$ws.Cells.Item(8, 4).Value = -1.392096468739544
$ws.Cells.Item(8, 5).Value = 0.1778806363245333
$ws.Cells.Item(8, 6).Value = 1.570796392423408
$ws.Cells.Item(9, 2).Value = 1.382027849918183
$ws.Cells.Item(9, 3).Value = 0.001143083554970447
$ws.Cells.Item(9, 4).Value = -1.389690162800365
$ws.Cells.Item(9, 5).Value = 0.179962940918031
$ws.Cells.Item(9, 6).Value = 1.570796392985272
$ws.Cells.Item(10, 2).Value = 1.472870163702626
$ws.Cells.Item(10, 3).Value = 0.001467079033438478
$ws.Cells.Item(10, 4).Value = -1.387283856861185
$ws.Cells.Item(10, 5).Value = 0.1820452455115286
$ws.Cells.Item(10, 6).Value = 1.570796393547136
$ws.Cells.Item(11, 2).Value = 1.558109815172848
$ws.Cells.Item(11, 3).Value = 0.00177109222412046
$ws.Cells.Item(11, 4).Value = -1.385025958850429
$ws.Cells.Item(11, 5).Value = 0.183999124826247
$ws.Cells.Item(11, 6).Value = 1.570796394074348
$ws.Cells.Item(12, 2).Value = 1.632681997596789
$ws.Cells.Item(12, 3).Value = 0.002037059138857806
$ws.Cells.Item(12, 4).Value = -1.38305062953539
$ws.Cells.Item(12, 5).Value = 0.1857084824101697
$ws.Cells.Item(12, 6).Value = 1.570796394535581
$ws.Cells.Item(13, 2).Value = 1.692597615406726
$ws.Cells.Item(13, 3).Value = 0.002250752388746849
$ws.Cells.Item(13, 4).Value = -1.381463535361107
$ws.Cells.Item(13, 5).Value = 0.1870818794648058
$ws.Cells.Item(13, 6).Value = 1.570796394906162
$ws.Cells.Item(14, 2).Value = 1.735481139781426
$ws.Cells.Item(14, 3).Value = 0.002403699483766303
$ws.Cells.Item(14, 4).Value = -1.380327601289231
$ws.Cells.Item(14, 5).Value = 0.1880648636719529
$ws.Cells.Item(14, 6).Value = 1.570796395171398
$ws.Cells.Item(15, 2).Value = 1.761108464228321
$ws.Cells.Item(15, 3).Value = 0.002495101132404725
$ws.Cells.Item(15, 4).Value = -1.379648763636899
$ws.Cells.Item(15, 5).Value = 0.1886522980204598
$ws.Cells.Item(15, 6).Value = 1.570796395329905
$ws.Cells.Item(16, 2).Value = 1.771944760165668
$ws.Cells.Item(16, 3).Value = 0.002533749541287982
$ws.Cells.Item(16, 4).Value = -1.379361722915605
$ws.Cells.Item(16, 5).Value = 0.1889006896329894
$ws.Cells.Item(16, 6).Value = 1.570796395396928
$ws.Cells.Item(17, 2).Value = 1.77368233250472
$ws.Cells.Item(17, 3).Value = 0.002539946714806697
$ws.Cells.Item(17, 4).Value = -1.379315696670069
$ws.Cells.Item(17, 5).Value = 0.1889405185927815
$ws.Cells.Item(17, 6).Value = 1.570796395407675
